$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.963.79'
$ws.Range("E2").Value = '  +10.55%  '
$ws.Range("D3").Value = '1.812.20'
$ws.Range("E3").Value = '  +7.40%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''228.00'
$ws.Range("E5").Value = '  +3.35%  '
$ws.Range("D6").Value = '''0.540'
$ws.Range("E6").Value = '  +3.50%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''30.85'
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("D9").Value = '''47.71'
$ws.Range("E9").Value = '  +7.71%  '
$ws.Range("E10").Value = '  +5.33%  '
$ws.Range("E11").Value = '  +6.60%  '
$ws.Range("D13").Value = '2.072.27'
$ws.Range("E13").Value = '  +7.33%  '
$ws.Range("D14").Value = '1.810.72'
$ws.Range("E14").Value = '  +7.34%  '
$ws.Range("D15").Value = '''0.637'
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("D16").Value = '34.003.04'
$ws.Range("E16").Value = '  +10.57%  '
$ws.Range("D17").Value = '''10.13'
$ws.Range("E17").Value = '  -4.14%  '
$ws.Range("D18").Value = '''4.26'
$ws.Range("E18").Value = '  +6.87%  '
$ws.Range("D19").Value = '''69.07'
$ws.Range("E19").Value = '  +3.90%  '
$ws.Range("D20").Value = '''255.51'
$ws.Range("E20").Value = '  +3.68%  '
$ws.Range("E21").Value = '  +3.93%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '''10.37'
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("D24").Value = '''4.30'
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("D26").Value = '''159.14'
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").Value = '''16.54'
$ws.Range("E27").Value = '  +4.32%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '''7.05'
$ws.Range("E28").Value = '  +5.12%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '''0.115'
$ws.Range("E29").Value = '  +3.26%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '''3.80'
$ws.Range("E31").Value = '  +9.08%  '
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("E33").Value = '  +5.65%  '
$ws.Range("E34").Value = '  +7.10%  '
$ws.Range("D35").Value = '1.546.67'
$ws.Range("E35").Value = '  +2.23%  '
$ws.Range("E36").Value = '  +4.02%  '
$ws.Range("E37").Value = '  +3.69%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.0186'
$ws.Range("E38").Value = '  +4.00%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").Value = '''83.92'
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("E40").Value = '  +5.75%  '
$ws.Range("D41").Value = '''2.84'
$ws.Range("E41").Value = '  +3.74%  '
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("D43").Value = '''0.901'
$ws.Range("E43").Value = '  +6.31%  '
$ws.Range("E44").Value = '  +4.97%  '
$ws.Range("D45").Value = '''0.0525'
$ws.Range("E45").Value = '  +4.87%  '
$ws.Range("E46").Value = '  +3.72%  '
$ws.Range("E47").Value = '  +7.69%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  +3.94%  '
$ws.Range("D50").Value = '''52.31'
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("E51").Value = '  +7.95%  '
